# "Resolved the Excel Issues"
# The "Users" sample data sheet is trimmed down to a single generic
# "string" placeholder row and the mailto: hyperlink on the old sample
# data is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old sample rows 3-5 (Jack Sparrow / Steven Cook / extra
# Software-Engineering row) - only the header row and one data row remain.
$ws.Rows("3:5").Delete()

# Drop the mailto: hyperlink that lived on C2 (Rajeev Singh's e-mail).
$ws.Hyperlinks.Delete()

# Replace the remaining data row with generic placeholder text and reset
# its formatting back to the workbook default ("Normal" - no explicit
# font/number-format overrides), same as the header cells lost their
# special hyperlink/date styling.
$ws.Range("A2:F2").Value = "string"
$ws.Range("A2:F2").Style = "Normal"
$ws.Rows(2).RowHeight = 13

# The "Hyperlink" cell style is no longer used anywhere once the
# hyperlinked cell is gone - drop it from the workbook's style gallery.
$wb.Styles.Item("Hyperlink").Delete()

# Move the active selection, matching the author's last cursor position.
$ws.Range("C8").Select() | Out-Null
